# Updates LR-pair (Fn1-Sdc2) data table to include full 4x4 sending/target
# cluster combinations per Dr Hou's advice (adds M2 as a target cluster and
# sCs as a sending cluster), refreshing all computed statistic columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.057109
$ws.Range("N2").Value = 9.171327
$ws.Range("O2").Value = 0.02694952608666365
$ws.Range("P2").Value = 0.02694952608666365
$ws.Range("Q2").Value = 61.00700279231067
$ws.Range("R2").Value = 549.063025130796
$ws.Range("S2").Value = 0.0003163158034925219
$ws.Range("T2").Value = 0.0003163158034925219

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 89.02756
$ws.Range("N3").Value = 267.08268
$ws.Range("O3").Value = 0.7848102735793893
$ws.Range("P3").Value = 0.7848102735793893
$ws.Range("Q3").Value = 1776.614638703627
$ws.Range("R3").Value = 15989.53174833264
$ws.Range("S3").Value = 0.009211586559189975
$ws.Range("T3").Value = 0.009211586559189975

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.184005
$ws.Range("N4").Value = 0.5520149999999999
$ws.Range("O4").Value = 0.0016220709001794
$ws.Range("P4").Value = 0.0016220709001794
$ws.Range("Q4").Value = 3.671963789579999
$ws.Range("R4").Value = 33.04767410621999
$ws.Range("S4").Value = 0.00001903880084800427
$ws.Range("T4").Value = 0.00001903880084800427

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 21.16964733333333
$ws.Range("N5").Value = 63.508942
$ws.Range("O5").Value = 0.1866181294337677
$ws.Range("P5").Value = 0.1866181294337677
$ws.Range("Q5").Value = 422.4568813139796
$ws.Range("R5").Value = 3802.111931825816
$ws.Range("S5").Value = 0.002190400802161996
$ws.Range("T5").Value = 0.002190400802161996

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.057109
$ws.Range("N6").Value = 9.171327
$ws.Range("O6").Value = 0.02694952608666365
$ws.Range("P6").Value = 0.02694952608666365
$ws.Range("Q6").Value = 5005.537070994423
$ws.Range("R6").Value = 45049.83363894981
$ws.Range("S6").Value = 0.02595325795488463
$ws.Range("T6").Value = 0.02595325795488463

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 89.02756
$ws.Range("N7").Value = 267.08268
$ws.Range("O7").Value = 0.7848102735793893
$ws.Range("P7").Value = 0.7848102735793893
$ws.Range("Q7").Value = 145768.6827392089
$ws.Range("R7").Value = 1311918.14465288
$ws.Range("S7").Value = 0.7557974641316252
$ws.Range("T7").Value = 0.7557974641316253

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Sdc2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.184005
$ws.Range("N8").Value = 0.5520149999999999
$ws.Range("O8").Value = 0.0016220709001794
$ws.Range("P8").Value = 0.0016220709001794
$ws.Range("Q8").Value = 301.2793618900499
$ws.Range("R8").Value = 2711.514257010449
$ws.Range("S8").Value = 0.001562106300425842
$ws.Range("T8").Value = 0.001562106300425842

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Sdc2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 21.16964733333333
$ws.Range("N9").Value = 63.508942
$ws.Range("O9").Value = 0.1866181294337677
$ws.Range("P9").Value = 0.1866181294337677
$ws.Range("Q9").Value = 34661.98114194758
$ws.Range("R9").Value = 311957.8302775283
$ws.Range("S9").Value = 0.1797192439183345
$ws.Range("T9").Value = 0.1797192439183345

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Sdc2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.057109
$ws.Range("N10").Value = 9.171327
$ws.Range("O10").Value = 0.02694952608666365
$ws.Range("P10").Value = 0.02694952608666365
$ws.Range("Q10").Value = 53.50191229130734
$ws.Range("R10").Value = 481.517210621766
$ws.Range("S10").Value = 0.0002774025865919831
$ws.Range("T10").Value = 0.0002774025865919832

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fn1"
$ws.Range("C11").Value = "Sdc2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 89.02756
$ws.Range("N11").Value = 267.08268
$ws.Range("O11").Value = 0.7848102735793893
$ws.Range("P11").Value = 0.7848102735793893
$ws.Range("Q11").Value = 1558.055243247494
$ws.Range("R11").Value = 14022.49718922744
$ws.Range("S11").Value = 0.008078375819106538
$ws.Range("T11").Value = 0.00807837581910654

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fn1"
$ws.Range("C12").Value = "Sdc2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.184005
$ws.Range("N12").Value = 0.5520149999999999
$ws.Range("O12").Value = 0.0016220709001794
$ws.Range("P12").Value = 0.0016220709001794
$ws.Range("Q12").Value = 3.22023826143
$ws.Range("R12").Value = 28.98214435287
$ws.Range("S12").Value = 0.0000166966447535426
$ws.Range("T12").Value = 0.0000166966447535426

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fn1"
$ws.Range("C13").Value = "Sdc2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 21.16964733333333
$ws.Range("N13").Value = 63.508942
$ws.Range("O13").Value = 0.1866181294337677
$ws.Range("P13").Value = 0.1866181294337677
$ws.Range("Q13").Value = 370.4861733310485
$ws.Range("R13").Value = 3334.375559979436
$ws.Range("S13").Value = 0.001920937371715154
$ws.Range("T13").Value = 0.001920937371715154

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fn1"
$ws.Range("C14").Value = "Sdc2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.057109
$ws.Range("N14").Value = 9.171327
$ws.Range("O14").Value = 0.02694952608666365
$ws.Range("P14").Value = 0.02694952608666365
$ws.Range("Q14").Value = 77.63871720744334
$ws.Range("R14").Value = 698.74845486699
$ws.Range("S14").Value = 0.0004025497416945135
$ws.Range("T14").Value = 0.0004025497416945135

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fn1"
$ws.Range("C15").Value = "Sdc2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 89.02756
$ws.Range("N15").Value = 267.08268
$ws.Range("O15").Value = 0.7848102735793893
$ws.Range("P15").Value = 0.7848102735793893
$ws.Range("Q15").Value = 2260.954893825733
$ws.Range("R15").Value = 20348.5940444316
$ws.Range("S15").Value = 0.01172284706946753
$ws.Range("T15").Value = 0.01172284706946753

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fn1"
$ws.Range("C16").Value = "Sdc2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.184005
$ws.Range("N16").Value = 0.5520149999999999
$ws.Range("O16").Value = 0.0016220709001794
$ws.Range("P16").Value = 0.0016220709001794
$ws.Range("Q16").Value = 4.673013673949999
$ws.Range("R16").Value = 42.05712306554999
$ws.Range("S16").Value = 0.0000242291541520106
$ws.Range("T16").Value = 0.00002422915415201059

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fn1"
$ws.Range("C17").Value = "Sdc2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 21.16964733333333
$ws.Range("N17").Value = 63.508942
$ws.Range("O17").Value = 0.1866181294337677
$ws.Range("P17").Value = 0.1866181294337677
$ws.Range("Q17").Value = 537.6269746005045
$ws.Range("R17").Value = 4838.642771404539
$ws.Range("S17").Value = 0.002787547341556118
$ws.Range("T17").Value = 0.002787547341556117

